$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$lo = $ws.ListObjects.Item(1)

# Insert a new row into the table (just above the totals row); this pushes
# the totals row down one row and grows the table range by one row.
$newListRow = $lo.ListRows.Add()

# Populate the new data row (row 8) with the added task.
$ws.Range("A8").Value = 7
$ws.Range("B8").Value = "simulate call spsavebaneditlog"
$ws.Range("C8").Value = 42993
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = "Vista"

# Re-create the totals row's SUM formula, now spanning the extra data row.
$ws.Range("D9").Formula = "=SUM(D2:D8)"

# Carry the date number-format over to the (now empty) totals row's "When"
# cell, matching the formatting the totals row used to carry in C8.
$ws.Range("C8").Copy()
$null = $ws.Range("C9").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Restore the last-used selection recorded in the workbook.
$null = $ws.Range("E15").Select()
